$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2004.85
$ws.Range("I70").Value = 1599
$ws.Range("J70").Value = 2049.9443
$ws.Range("K70").Value = 4797
$ws.Range("L70").Value = 6149.8329
$ws.Range("M70").Value = -4527
$ws.Range("N70").Value = -6689.8329

$ws.Range("H73").Value = 2004.85
$ws.Range("I73").Value = 1599
$ws.Range("J73").Value = 2049.9443
$ws.Range("K73").Value = 4797
$ws.Range("L73").Value = 6149.8329
$ws.Range("M73").Value = -3861
$ws.Range("N73").Value = -8021.8329

$ws.Range("H132").Value = 1492
$ws.Range("I132").Value = 1491.125
$ws.Range("J132").Value = 1499
$ws.Range("K132").Value = 4473.375
$ws.Range("L132").Value = 4497
$ws.Range("M132").Value = -1943.375
$ws.Range("N132").Value = -9557

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1395.7273
$ws.Range("I61").Value = 1426.7778
$ws.Range("J61").Value = 1256
$ws.Range("K61").Value = 1426.7778
$ws.Range("L61").Value = 1256
$ws.Range("M61").Value = -1214.7778
$ws.Range("N61").Value = -1680

$ws.Range("H136").Value = 1395.7273
$ws.Range("I136").Value = 1426.7778
$ws.Range("J136").Value = 1256
$ws.Range("K136").Value = 4280.3334
$ws.Range("L136").Value = 3768
$ws.Range("M136").Value = -1730.3334
$ws.Range("N136").Value = -8868

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4387.25
$ws.Range("J86").Value = 5599.6665
$ws.Range("L86").Value = 5599.6665
$ws.Range("N86").Value = -7845.6665

$ws.Range("H89").Value = 4387.25
$ws.Range("J89").Value = 5599.6665
$ws.Range("L89").Value = 27998.3325
$ws.Range("N89").Value = -39230.3325

$ws.Range("H99").Value = 1144.3572
$ws.Range("I99").Value = 1162.4615
$ws.Range("K99").Value = 1162.4615
$ws.Range("M99").Value = 335.5385000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1413.15
$ws.Range("I7").Value = 753.9375
$ws.Range("K7").Value = 753.9375
$ws.Range("M7").Value = -640.9375

$ws.Range("H22").Value = 1730.7693
$ws.Range("I22").Value = 505.55554
$ws.Range("K22").Value = 505.55554
$ws.Range("M22").Value = -155.55554

$ws.Range("H62").Value = 1834.6666
$ws.Range("I62").Value = 2002.5
$ws.Range("J62").Value = 1499
$ws.Range("K62").Value = 2002.5
$ws.Range("L62").Value = 1499
$ws.Range("M62").Value = -1378.5
$ws.Range("N62").Value = -2747

$ws.Range("H65").Value = 1834.6666
$ws.Range("I65").Value = 2002.5
$ws.Range("J65").Value = 1499
$ws.Range("K65").Value = 10012.5
$ws.Range("L65").Value = 7495
$ws.Range("M65").Value = -6892.5
$ws.Range("N65").Value = -13735

$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws.Range("H99").Value = 4132.222
$ws.Range("I99").Value = 4000
$ws.Range("J99").Value = 4238
$ws.Range("K99").Value = 4000
$ws.Range("L99").Value = 4238
$ws.Range("M99").Value = -2502
$ws.Range("N99").Value = -7234

$ws.Range("H122").Value = 1103
$ws.Range("I122").Value = 637.3333
$ws.Range("K122").Value = 1911.9999
$ws.Range("M122").Value = 538.0001

$ws.Range("H126").Value = 4132.222
$ws.Range("I126").Value = 4000
$ws.Range("J126").Value = 4238
$ws.Range("K126").Value = 12000
$ws.Range("L126").Value = 12714
$ws.Range("M126").Value = -9530
$ws.Range("N126").Value = -17654

$ws.Range("H132").Value = 1504.4
$ws.Range("I132").Value = 1504.4
$ws.Range("K132").Value = 4513.200000000001
$ws.Range("M132").Value = -1983.200000000001

$ws.Range("H134").Value = 1200
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1397.1
$ws.Range("I23").Value = 1374.5
$ws.Range("J23").Value = 1412.1666
$ws.Range("K23").Value = 4123.5
$ws.Range("L23").Value = 4236.4998
$ws.Range("M23").Value = -3888.5
$ws.Range("N23").Value = -4706.4998

$ws.Range("H33").Value = 198.08333
$ws.Range("I33").Value = 156.6
$ws.Range("J33").Value = 227.71428
$ws.Range("K33").Value = 939.5999999999999
$ws.Range("L33").Value = 1366.28568
$ws.Range("M33").Value = -656.5999999999999
$ws.Range("N33").Value = -1932.28568

$ws.Range("H68").Value = 3670.9285
$ws.Range("J68").Value = 3837.8
$ws.Range("L68").Value = 11513.4
$ws.Range("N68").Value = -13135.4

$ws.Range("H71").Value = 3670.9285
$ws.Range("J71").Value = 3837.8
$ws.Range("L71").Value = 34540.2
$ws.Range("N71").Value = -42652.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 2000
$ws.Range("I31").Value = 2000
$ws.Range("K31").Value = 2000
$ws.Range("M31").Value = -1708

$ws.Range("H37").Value = 2000
$ws.Range("I37").Value = 2000
$ws.Range("K37").Value = 2000
$ws.Range("M37").Value = -1723

$ws.Range("H80").Value = 6417.6665
$ws.Range("I80").Value = 1500
$ws.Range("J80").Value = 7401.2
$ws.Range("K80").Value = 1500
$ws.Range("L80").Value = 7401.2
$ws.Range("M80").Value = -502
$ws.Range("N80").Value = -9397.200000000001

$ws.Range("H83").Value = 6417.6665
$ws.Range("I83").Value = 1500
$ws.Range("J83").Value = 7401.2
$ws.Range("K83").Value = 7500
$ws.Range("L83").Value = 37006
$ws.Range("M83").Value = -2508
$ws.Range("N83").Value = -46990

$ws.Range("H132").Value = 4304.1665
$ws.Range("I132").Value = 4471.4443
$ws.Range("K132").Value = 13414.3329
$ws.Range("M132").Value = -10884.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4147.1
$ws.Range("I40").Value = 4163.4443
$ws.Range("K40").Value = 4163.4443
$ws.Range("M40").Value = -4027.4443

$ws.Range("H68").Value = 3600.6
$ws.Range("J68").Value = 3875.75
$ws.Range("L68").Value = 3875.75
$ws.Range("N68").Value = -5373.75

$ws.Range("H71").Value = 3600.6
$ws.Range("J71").Value = 3875.75
$ws.Range("L71").Value = 19378.75
$ws.Range("N71").Value = -26866.75

$ws.Range("H122").Value = 6177.778
$ws.Range("J122").Value = 6420
$ws.Range("L122").Value = 19260
$ws.Range("N122").Value = -24160

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

$ws.Range("H122").Value = 1703.6842
$ws.Range("I122").Value = 1830.5
$ws.Range("K122").Value = 5491.5
$ws.Range("M122").Value = -3041.5

$ws.Range("H132").Value = 3333.3333
$ws.Range("J132").Value = 3666.6667
$ws.Range("L132").Value = 11000.0001
$ws.Range("N132").Value = -16060.0001

$ws.Range("H135").Value = 65000
$ws.Range("L135").Value = 65000
$ws.Range("N135").Value = -75140
